# Recording protocol for 3 channels: add a new "Sequence5_singleDOF" column
# (column E) to the "Sequences" worksheet, mirroring the single-DOF
# Sequence1 pattern but using only the "Fist" DOF (Rest / Fist alternating).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sequences")

$ws.Range("E1").Value = "Sequence5_singleDOF"
$ws.Range("E2").Value = "Rest"
$ws.Range("E3").Value = "Fist"
$ws.Range("E4").Value = "Rest"
$ws.Range("E5").Value = "Fist"
$ws.Range("E6").Value = "Rest"
$ws.Range("E7").Value = "Fist"
$ws.Range("E8").Value = "Rest"
$ws.Range("E9").Value = "Fist"
$ws.Range("E10").Value = "Rest"

# Leave the selection where the author left it.
[void]$ws.Range("E13").Select()
